$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 5, shifting existing rows 5-36 down to 6-37.
$ws.Rows.Item(5).Insert()

# Populate the newly inserted row 5 with the new weekly record.
$ws.Range("A5").Value = 10
$ws.Range("B5").Value = "Vega Modelo de Temuco"
$ws.Range("C5").Value = "La Araucanía"
$ws.Range("D5").Value = 44490
$ws.Range("E5").Value = 9
$ws.Range("F5").Value = 100112026
$ws.Range("G5").Value = "Haba"
$ws.Range("H5").Value = "Sin especificar"
$ws.Range("I5").Value = "Primera"
$ws.Range("J5").Value = 65
$ws.Range("K5").Value = 9000
$ws.Range("L5").Value = 9000
$ws.Range("M5").Value = 9000
$ws.Range("N5").Value = "$/saco 25 kilos"
$ws.Range("O5").Value = "Región Metropolitana"
$ws.Range("P5").Value = 360
$ws.Range("Q5").Value = 25
$ws.Range("R5").Value = "Hortaliza"

# Apply the same date-number format the rest of column D uses.
$ws.Range("D5").NumberFormat = $ws.Range("D6").NumberFormat
